# energy_calc.xlsx update — "updates to IEAGHG scenario"
#
# 1. Insert a new "Heat Flare" sheet right after "Heat Recovery" (becomes
#    the 2nd tab), modeled on the "bb heat" layout, with its own
#    KnownQty/UnknownQty row describing a heat -> waste heat flow.
# 2. Select the full header row on "Heat Recovery" (row 1).
# 3. Leave the "Heat Flare" sheet as the active tab with F2 selected.

$wb = $excel.ActiveWorkbook

# --- Heat Recovery: select the entire first row ---------------------------
$heatRecovery = $wb.Worksheets.Item(1)
$heatRecovery.Rows.Item(1).Select() | Out-Null

# --- New sheet: copy "bb heat" (3rd tab) to just before "bb electricity" --
# (2nd tab) so it lands as the 2nd sheet overall, then rename it.
$bbHeat = $wb.Worksheets.Item(3)
$bbElectricity = $wb.Worksheets.Item(2)
$bbHeat.Copy($bbElectricity)

$heatFlare = $wb.Worksheets.Item(2)
$heatFlare.Name = "Heat Flare"

# Replace the copied data row with the new Heat Flare values.
$heatFlare.Rows.Item(2).Delete()
$heatFlare.Cells.Item(2, 1).Value = "heat"
$heatFlare.Cells.Item(2, 2).Value = "inflows"
$heatFlare.Cells.Item(2, 3).Value = "waste heat"
$heatFlare.Cells.Item(2, 4).Value = "outflows"
$heatFlare.Cells.Item(2, 5).Value = "returnvalue"
$heatFlare.Cells.Item(2, 6).Value = "none"

$heatFlare.Range("F2").Select() | Out-Null
